$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the minor-maj7 chords first (still at rows 63/64 at this point,
# before the "extended open minor chords" block gets shifted down) so the
# new shared-string entries land in the same order the original edit did.
$ws.Range("A63").Value = "Am(maj7)"
$ws.Range("A64").Value = "Em(maj7)"

# Rename the existing maj7 chords to use the "(maj7)" notation
$ws.Range("A51").Value = "A(maj7)"
$ws.Range("A52").Value = "D(maj7)"
$ws.Range("A53").Value = "G(maj7)"
$ws.Range("A55").Value = "E(maj7)"
$ws.Range("A56").Value = "C(maj7)"

# Insert two new rows at 58:59 (leaving the existing blank row 57 alone).
# This shifts the "extended open minor chords" block (rows 58-64) down to
# rows 60-66, making room for a new B7 row at 57 and leaving a two-row gap
# at 58-59 (matching the original sheet's blank-row separators).
$ws.Rows("58:59").Insert()

# New row 57: B7 chord
$ws.Range("A57").Value = "B7"
$ws.Range("B57").Value = "021202"
$ws.Range("B57").NumberFormat = "@"

# New row 67: Bm7 chord, appended at the end
$ws.Range("A67").Value = "Bm7"
$ws.Range("B67").Value = "020202"
$ws.Range("B67").NumberFormat = "@"

# Touch UsedRange to flush any stray formatting left behind in the blank
# rows created by Insert() above.
$null = $ws.UsedRange.Rows.Count

# Update the saved selection to match where editing left off
$ws.Range("C65").Select() | Out-Null
